$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("VT-SaleVoid-DualCF-Generic")
$ws.Cells.Item(2, 2).Value = "Thu Sep 04 07:41:01 IST 2025"
$ws.Cells.Item(3, 2).Value = "Thu Sep 04 07:42:00 IST 2025"
$ws.Cells.Item(4, 2).Value = "Thu Sep 04 07:42:51 IST 2025"
$ws.Cells.Item(5, 2).Value = "Thu Sep 04 07:43:49 IST 2025"
$ws = $wb.Worksheets.Item("VT-SaleVoid-NoCF-Generic")
$ws.Cells.Item(2, 2).Value = "Thu Sep 04 07:44:42 IST 2025"
$ws.Cells.Item(3, 2).Value = "Thu Sep 04 07:45:51 IST 2025"
$ws.Cells.Item(4, 2).Value = "Thu Sep 04 07:46:52 IST 2025"
$ws.Cells.Item(5, 2).Value = "Thu Sep 04 07:47:57 IST 2025"
$ws = $wb.Worksheets.Item("VT-SaleVoid-SingleCF-Generic")
$ws.Cells.Item(2, 2).Value = "Thu Sep 04 07:48:53 IST 2025"
$ws.Cells.Item(3, 2).Value = "Thu Sep 04 07:49:42 IST 2025"
$ws.Cells.Item(4, 2).Value = "Thu Sep 04 07:50:52 IST 2025"
$ws.Cells.Item(5, 2).Value = "Thu Sep 04 07:51:57 IST 2025"
$ws = $wb.Worksheets.Item("VT-SaleCredit-DualCF-Generic")
$ws.Cells.Item(2, 2).Value = "Thu Sep 04 07:33:44 IST 2025"
$ws.Cells.Item(3, 2).Value = "Thu Sep 04 07:34:39 IST 2025"
$ws.Cells.Item(4, 2).Value = "Thu Sep 04 07:35:39 IST 2025"
$ws.Cells.Item(5, 2).Value = "Thu Sep 04 07:36:35 IST 2025"
$ws = $wb.Worksheets.Item("VT-SaleCredit-SingleCF-Generic")
$ws.Cells.Item(2, 2).Value = "Thu Sep 04 07:37:28 IST 2025"
$ws.Cells.Item(3, 2).Value = "Thu Sep 04 07:38:24 IST 2025"
$ws.Cells.Item(4, 2).Value = "Thu Sep 04 07:39:21 IST 2025"
$ws.Cells.Item(5, 2).Value = "Thu Sep 04 07:40:10 IST 2025"
$ws = $wb.Worksheets.Item("VT-AuthCapCredit-Generic")
$ws.Cells.Item(2, 2).Value = "Thu Sep 04 06:41:45 IST 2025"
$ws.Cells.Item(3, 2).Value = "Thu Sep 04 06:43:13 IST 2025"
$ws.Cells.Item(4, 2).Value = "Thu Sep 04 06:44:36 IST 2025"
$ws.Cells.Item(5, 2).Value = "Thu Sep 04 06:46:00 IST 2025"
$ws.Cells.Item(6, 2).Value = "Thu Sep 04 06:47:20 IST 2025"
$ws.Cells.Item(7, 2).Value = "Thu Sep 04 06:48:34 IST 2025"
$ws = $wb.Worksheets.Item("VT-AuthCapVoid-Generic")
$ws.Cells.Item(2, 2).Value = "Thu Sep 04 06:49:49 IST 2025"
$ws.Cells.Item(3, 2).Value = "Thu Sep 04 06:51:09 IST 2025"
$ws.Cells.Item(4, 2).Value = "Thu Sep 04 06:52:18 IST 2025"
$ws.Cells.Item(5, 2).Value = "Thu Sep 04 06:53:29 IST 2025"
$ws.Cells.Item(6, 2).Value = "Thu Sep 04 06:54:45 IST 2025"
$ws.Cells.Item(7, 2).Value = "Thu Sep 04 06:56:00 IST 2025"
$ws = $wb.Worksheets.Item("VT-ManualAuthCapture-Generic")
$ws.Cells.Item(2, 2).Value = "Thu Sep 04 07:27:51 IST 2025"
$ws.Cells.Item(3, 2).Value = "Thu Sep 04 07:28:47 IST 2025"
$ws.Cells.Item(4, 2).Value = "Thu Sep 04 07:29:38 IST 2025"
$ws.Cells.Item(5, 2).Value = "Thu Sep 04 07:30:49 IST 2025"
$ws.Cells.Item(6, 2).Value = "Thu Sep 04 07:31:44 IST 2025"
$ws.Cells.Item(7, 2).Value = "Thu Sep 04 07:32:47 IST 2025"

$wb.Save()
